$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new header cells (Wins, Losses, Ties) in columns AD, AE, AF of row 1
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Copy the formatting of the existing header cell (AC1) onto the new headers
# so they match the bold/centered/bordered style used by the rest of row 1.
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

# Fill in the season record (Wins, Losses, Ties) for every player row (2-49)
for ($r = 2; $r -le 49; $r++) {
    $ws.Cells.Item($r, 30).Value = 87
    $ws.Cells.Item($r, 31).Value = 75
    $ws.Cells.Item($r, 32).Value = 0
}
